$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression (only B2 value changes slightly)
$ws.Range("B2").Value = 0.08559727483388412

# Row 3: RandomForestRegressor
$ws.Range("B3").Value = 0.02253158136884355
$ws.Range("C3").Value = 0.02284867209292149
$ws.Range("D3").Value = 0.02448936403410701

# Row 4: model name changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02346950018689528
$ws.Range("C4").Value = 0.0228719925136173
$ws.Range("D4").Value = 0.02290051981985665

# Row 5: model name changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01957789447367292
$ws.Range("C5").Value = 0.0186877691308285
$ws.Range("D5").Value = 0.01734001657310653
